# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型
# sheets to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1547
$ws1.Range("F6").Value  = 103
$ws1.Range("F8").Value  = 6343
$ws1.Range("F10").Value = 411
$ws1.Range("F12").Value = 5444
$ws1.Range("F15").Value = 1205
$ws1.Range("F24").Value = 3896

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1547
$ws4.Range("F7").Value  = 103
$ws4.Range("F9").Value  = 6343
$ws4.Range("F11").Value = 411
$ws4.Range("F13").Value = 5444
$ws4.Range("F16").Value = 1205
$ws4.Range("F25").Value = 3896
